$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-23 down to 5-24
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new eLibrary id entry for Бобов Д.Г.
$ws.Cells.Item(4, 1).Value = 6187803222
$ws.Cells.Item(4, 2).Value = "Бобов Д.Г., Бобов Дмитрий Геннадиевич, Бобов Дмитрий Геннадьевич"
